# Refresh the crypto price / 1h-volume snapshot (GitHub Actions bot update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Cells.Item(2, 4).Value = "72.146.17"
$ws.Cells.Item(2, 5).Value = "  +3.83%  "

# Row 3: Ethereum
$ws.Cells.Item(3, 4).Value = "2.616.62"
$ws.Cells.Item(3, 5).Value = "  +4.07%  "

# Row 4: TetherUSD
$ws.Cells.Item(4, 5).Value = "  -0.02%  "

# Row 5: BNB
$ws.Cells.Item(5, 4).Value = "'603.25"
$ws.Cells.Item(5, 5).Value = "  +0.92%  "

# Row 6: Solana
$ws.Cells.Item(6, 4).Value = "'178.00"
$ws.Cells.Item(6, 5).Value = "  +0.88%  "

# Row 7: USDC
$ws.Cells.Item(7, 5).Value = "  +0.02%  "

# Row 8: XRP
$ws.Cells.Item(8, 4).Value = "'0.524"
$ws.Cells.Item(8, 5).Value = "  +1.08%  "

# Row 9: now Dogecoin (swapped with row 10)
$ws.Cells.Item(9, 2).Value = "Dogecoin"
$ws.Cells.Item(9, 3).Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Cells.Item(9, 4).Value = "'0.171"
$ws.Cells.Item(9, 5).Value = "  +8.79%  "

# Row 10: now LidoStakedEther (swapped with row 9)
$ws.Cells.Item(10, 2).Value = "LidoStakedEther"
$ws.Cells.Item(10, 3).Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Cells.Item(10, 4).Value = "2.614.87"
$ws.Cells.Item(10, 5).Value = "  +4.01%  "

# Row 11: TRON
$ws.Cells.Item(11, 4).Value = "'0.166"
$ws.Cells.Item(11, 5).Value = "  +0.93%  "

# Row 12: Cardano
$ws.Cells.Item(12, 4).Value = "'0.351"
$ws.Cells.Item(12, 5).Value = "  +2.74%  "

# Row 13: Toncoin
$ws.Cells.Item(13, 4).Value = "'5.03"
$ws.Cells.Item(13, 5).Value = "  +0.52%  "

# Row 14: WrappedliquidstakedEther2.0
$ws.Cells.Item(14, 4).Value = "3.095.71"
$ws.Cells.Item(14, 5).Value = "  +4.19%  "

# Row 15: ShibaInu
$ws.Cells.Item(15, 4).Value = "'0.0000187"
$ws.Cells.Item(15, 5).Value = "  +4.59%  "

# Row 16: WrappedBTC
$ws.Cells.Item(16, 4).Value = "72.036.72"
$ws.Cells.Item(16, 5).Value = "  +3.63%  "

# Row 17: Avalanche
$ws.Cells.Item(17, 4).Value = "'26.46"
$ws.Cells.Item(17, 5).Value = "  +1.92%  "

# Row 18: WrappedEther
$ws.Cells.Item(18, 4).Value = "2.617.70"
$ws.Cells.Item(18, 5).Value = "  +5.35%  "

# Row 19: BitcoinCash
$ws.Cells.Item(19, 4).Value = "'382.03"
$ws.Cells.Item(19, 5).Value = "  +5.49%  "

# Row 20: Chainlink
$ws.Cells.Item(20, 4).Value = "'11.54"
$ws.Cells.Item(20, 5).Value = "  +4.77%  "

# Row 21: Uniswap
$ws.Cells.Item(21, 4).Value = "'7.89"
$ws.Cells.Item(21, 5).Value = "  +3.02%  "

# Row 22: Polkadot
$ws.Cells.Item(22, 4).Value = "'4.16"
$ws.Cells.Item(22, 5).Value = "  +1.84%  "

# Row 23: SuiNetwork
$ws.Cells.Item(23, 4).Value = "'1.99"
$ws.Cells.Item(23, 5).Value = "  +17.42%  "

# Row 24: Litecoin
$ws.Cells.Item(24, 4).Value = "'72.97"
$ws.Cells.Item(24, 5).Value = "  +3.55%  "

# Row 25: Dai
$ws.Cells.Item(25, 5).Value = "  -0.14%  "

# Row 26: NEARProtocol
$ws.Cells.Item(26, 4).Value = "'4.36"
$ws.Cells.Item(26, 5).Value = "  +3.18%  "

# Row 27: Aptos
$ws.Cells.Item(27, 4).Value = "'9.84"
$ws.Cells.Item(27, 5).Value = "  +8.71%  "

# Row 28: WrappedeETH
$ws.Cells.Item(28, 4).Value = "2.751.50"
$ws.Cells.Item(28, 5).Value = "  +5.57%  "

# Row 29: Binance-PegBSC-USD
$ws.Cells.Item(29, 5).Value = "  -0.03%  "

# Row 30: PEPE
$ws.Cells.Item(30, 4).Value = "0.0₃0946"
$ws.Cells.Item(30, 5).Value = "  +5.55%  "

# Row 31: Bittensor
$ws.Cells.Item(31, 4).Value = "'518.00"
$ws.Cells.Item(31, 5).Value = "  +1.07%  "

# Row 32: InternetComputer(DFINITY)
$ws.Cells.Item(32, 4).Value = "'8.02"
$ws.Cells.Item(32, 5).Value = "  +3.52%  "

# Row 33: Fetch.AI
$ws.Cells.Item(33, 5).Value = "  +6.22%  "

# Row 34: PancakeSwap
$ws.Cells.Item(34, 5).Value = "  +2.27%  "

# Row 35: FirstDigitalUSD
$ws.Cells.Item(35, 4).Value = "'0.998"
$ws.Cells.Item(35, 5).Value = "  -0.17%  "

# Row 36: Monero
$ws.Cells.Item(36, 4).Value = "'163.57"
$ws.Cells.Item(36, 5).Value = "  +1.05%  "

# Row 37: EthereumClassic
$ws.Cells.Item(37, 4).Value = "'19.24"
$ws.Cells.Item(37, 5).Value = "  +2.74%  "

# Row 38: WhiteBITCoin
$ws.Cells.Item(38, 4).Value = "'19.08"
$ws.Cells.Item(38, 5).Value = "  +1.01%  "

# Row 39: now ImmutableX (swapped with row 40)
$ws.Cells.Item(39, 2).Value = "ImmutableX"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(39, 4).Value = "'1.40"
$ws.Cells.Item(39, 5).Value = "  +5.88%  "

# Row 40: now Kaspa (swapped with row 39)
$ws.Cells.Item(40, 2).Value = "Kaspa"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(40, 4).Value = "'0.112"
$ws.Cells.Item(40, 5).Value = "  -6.15%  "

# Row 41: Stacks
$ws.Cells.Item(41, 4).Value = "'1.82"
$ws.Cells.Item(41, 5).Value = "  +6.04%  "

# Row 42: USDe
$ws.Cells.Item(42, 5).Value = "  -0.04%  "

# Row 43: RenderToken
$ws.Cells.Item(43, 4).Value = "'5.03"
$ws.Cells.Item(43, 5).Value = "  +4.77%  "

# Row 44: dogwifhat
$ws.Cells.Item(44, 4).Value = "'2.57"
$ws.Cells.Item(44, 5).Value = "  +8.42%  "

# Row 45: PolygonEcosystemToken
$ws.Cells.Item(45, 4).Value = "'0.331"
$ws.Cells.Item(45, 5).Value = "  +3.25%  "

# Row 46: OKB
$ws.Cells.Item(46, 4).Value = "'39.47"
$ws.Cells.Item(46, 5).Value = "  +1.70%  "

# Row 47: Aave
$ws.Cells.Item(47, 4).Value = "'150.22"
$ws.Cells.Item(47, 5).Value = "  +0.23%  "

# Row 48: Filecoin
$ws.Cells.Item(48, 5).Value = "  +2.59%  "

# Row 49: ARBITRUM
$ws.Cells.Item(49, 4).Value = "'0.540"
$ws.Cells.Item(49, 5).Value = "  +4.72%  "

# Row 50: Optimism
$ws.Cells.Item(50, 4).Value = "'1.69"
$ws.Cells.Item(50, 5).Value = "  +7.11%  "

# Row 51: Cronos
$ws.Cells.Item(51, 5).Value = "  +3.52%  "

Write-Host "Applied cryptos list update."
